$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells hold their numeric-looking values as TEXT (shared strings) in the
# source workbook, so force text formatting before assigning, otherwise
# Excel auto-converts the values to numbers. Restore the original
# (default/"Normal") cell style afterwards so formatting is unaffected.
$cells = @("B11","C11","D11","B12","C12","D12","B14","D14")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Enterprises density (per 1000 people) - row 11
$ws.Range("B11").Value = "23.19"
$ws.Range("C11").Value = "2.05"
$ws.Range("D11").Value = "25.25"

# Employment (% of total) - row 12
$ws.Range("B12").Value = "40.97"
$ws.Range("C12").Value = "24.71"
$ws.Range("D12").Value = "65.68"

# Enterprises (% of total) - row 14
$ws.Range("B14").Value = "91.47"
$ws.Range("D14").Value = "99.57"

foreach ($addr in $cells) {
    $ws.Range($addr).Style = "Normal"
}
